$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table")

# Row 4: new "Name in Files" entry, and "Concluded?" flipped from "Tem Only" to "Yes"
$ws.Range("A4").Value = "mcmc_fixshift_predictors/D_mcmc"
$ws.Range("P4").Value = "Yes"
$ws.Range("P4").Font.Bold = $false

# Row 6: rename restored B_covar_rjmcmc run, diagnostics now "Good"
$ws.Range("A6").Value = "mcmc_predictors/B_covar_rjmcmc"
$ws.Range("Q6").Value = "Good"

# Row 7: rename restored B_covar_rjmcmc_hpp run, diagnostics now "Good"
$ws.Range("A7").Value = "mcmc_predictors/B_covar_rjmcmc_hpp_notranslate"
$ws.Range("Q7").Value = "Good"
$ws.Rows.Item(7).RowHeight = 12.45

# Update the selected cell to match the saved view state
[void]$ws.Range("N25").Select()
